{"js": "// Update the \"Tempo total do processo\" day counters and renumber a batch of\n// contract reference numbers (785810/2024-0NN/00 -> 785810/2024-0MM/00),\n// as described in the commit \"atualizacao matriz de risco\".\n\nconst replacements = [\n  // Day-count updates inside the \"Tempo total do processo\" block\n  [\"Conclu\u00eddo (97 dias)\", \"Conclu\u00eddo (110 dias)\"],\n  [\"Assinatura Contrato (106 dias)\", \"Assinatura Contrato (119 dias)\"],\n  [\"Assinatura Contrato (14 dias)\", \"Assinatura Contrato (27 dias)\"],\n  [\"Conclu\u00eddo (22 dias)\", \"Conclu\u00eddo (35 dias)\"],\n  [\"Assinatura Contrato (49 dias)\", \"Assinatura Contrato (62 dias)\"],\n  [\"Total de dias 708\", \"Total de dias 773\"],\n\n  // Contract / process number renumbering in the \"Rela\u00e7\u00e3o de ...\" section\n  [\"785810/2024-010/00\", \"785810/2024-055/00\"],\n  [\"785810/2024-011/00\", \"785810/2024-056/00\"],\n  [\"785810/2024-012/00\", \"785810/2024-057/00\"],\n  [\"785810/2024-013/00\", \"785810/2024-058/00\"],\n  [\"785810/2024-014/00\", \"785810/2024-059/00\"],\n  [\"785810/2024-015/00\", \"785810/2024-060/00\"],\n  [\"785810/2024-016/00\", \"785810/2024-061/00\"],\n  [\"785810/2024-017/00\", \"785810/2024-062/00\"],\n  [\"785810/2024-018/00\", \"785810/2024-063/00\"],\n  [\"785810/2024-019/00\", \"785810/2024-064/00\"],\n  [\"785810/2024-020/00\", \"785810/2024-065/00\"],\n  [\"785810/2024-021/00\", \"785810/2024-066/00\"],\n];\n\nfor (const [searchText, newText] of replacements) {\n  const results = context.document.body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the \"Tempo total do processo\" day counters and renumber a batch of\n# contract reference numbers (785810/2024-0NN/00 -> 785810/2024-0MM/00),\n# as described in the commit \"atualizacao matriz de risco\".\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = @(\n    # Day-count updates inside the \"Tempo total do processo\" block\n    @(\"Conclu\u00eddo (97 dias)\", \"Conclu\u00eddo (110 dias)\"),\n    @(\"Assinatura Contrato (106 dias)\", \"Assinatura Contrato (119 dias)\"),\n    @(\"Assinatura Contrato (14 dias)\", \"Assinatura Contrato (27 dias)\"),\n    @(\"Conclu\u00eddo (22 dias)\", \"Conclu\u00eddo (35 dias)\"),\n    @(\"Assinatura Contrato (49 dias)\", \"Assinatura Contrato (62 dias)\"),\n    @(\"Total de dias 708\", \"Total de dias 773\"),\n\n    # Contract / process number renumbering in the \"Rela\u00e7\u00e3o de ...\" section\n    @(\"785810/2024-010/00\", \"785810/2024-055/00\"),\n    @(\"785810/2024-011/00\", \"785810/2024-056/00\"),\n    @(\"785810/2024-012/00\", \"785810/2024-057/00\"),\n    @(\"785810/2024-013/00\", \"785810/2024-058/00\"),\n    @(\"785810/2024-014/00\", \"785810/2024-059/00\"),\n    @(\"785810/2024-015/00\", \"785810/2024-060/00\"),\n    @(\"785810/2024-016/00\", \"785810/2024-061/00\"),\n    @(\"785810/2024-017/00\", \"785810/2024-062/00\"),\n    @(\"785810/2024-018/00\", \"785810/2024-063/00\"),\n    @(\"785810/2024-019/00\", \"785810/2024-064/00\"),\n    @(\"785810/2024-020/00\", \"785810/2024-065/00\"),\n    @(\"785810/2024-021/00\", \"785810/2024-066/00\")\n)\n\nforeach ($pair in $replacements) {\n    $searchText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    $find.Execute($searchText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll)\n}\n"}
